$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.859.40'
$ws.Range("E2").Value = '  -4.86%  '

$ws.Range("D3").Value = '1.953.38'
$ws.Range("E3").Value = '  -4.99%  '

$ws.Range("E4").Value = '  +0.18%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '241.46'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  -4.56%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '0.624'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -3.93%  '

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '62.10'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  -6.33%  '

$ws.Range("E8").Value = '  +0.11%  '

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.367'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  -2.63%  '

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '56.00'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  -5.81%  '

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.0796'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +5.26%  '

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '0.852'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  -6.65%  '

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '22.07'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +5.71%  '

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '13.96'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  -9.00%  '

$ws.Range("D16").Value = '2.242.54'
$ws.Range("E16").Value = '  -4.78%  '

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '5.40'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  -4.07%  '

$ws.Range("D18").Value = '1.970.90'
$ws.Range("E18").Value = '  -3.92%  '

$ws.Range("D19").Value = '35.718.20'
$ws.Range("E19").Value = '  -4.88%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '70.96'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -3.12%  '

$ws.Range("D21").Value = '0.0₃0851'
$ws.Range("E21").Value = '  -3.24%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '237.16'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  -0.28%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '5.18'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -3.52%  '

$ws.Range("E24").Value = '  -0.07%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '2.51'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  -9.77%  '

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '2.29'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  -2.72%  '

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '9.79'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  +1.72%  '

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '158.93'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  -3.97%  '

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '19.72'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -1.22%  '

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '0.128'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  +13.84%  '

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '0.119'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  -2.24%  '

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '4.84'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -7.71%  '

$ws.Range("E33").Value = '  -7.66%  '

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '0.0617'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +0.02%  '

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '4.39'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  -7.76%  '

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '6.24'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  +3.71%  '

$ws.Range("B37").Value = 'BinanceUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  +0.30%  '

$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '2.27'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  -7.69%  '

$ws.Range("E39").Value = '  +1.02%  '

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '3.11'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  +14.70%  '

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.0983'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  -5.76%  '

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '1.22'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -1.70%  '

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '0.0211'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -3.63%  '

$ws.Range("E44").Value = '  -4.56%  '

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '1.08'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  -5.30%  '

$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '91.86'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -4.06%  '

$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '16.09'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  -6.16%  '

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '7.53'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  -7.94%  '

$ws.Range("D49").Value = '1.334.95'
$ws.Range("E49").Value = '  -6.74%  '

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '2.76'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -6.09%  '

$ws.Range("D51").Value = '2.138.09'
$ws.Range("E51").Value = '  -4.62%  '
